# Daily attendance processing - move leading "System" entry in the
# "Recorded By" column (G) to the end of the comma-separated list.
#
# Example: "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $value = $cell.Text

    if ($value -and $value.StartsWith("System, ")) {
        $parts = $value -split ", "
        $reordered = ($parts[1..($parts.Length - 1)] + $parts[0]) -join ", "
        $cell.Value = $reordered
    }
}
